$d = $word.ActiveDocument

# Helper: Word represents a manual line break (<w:br/>) as the caret code ^l
# both in Find What and Replace With strings.

function Replace-InParagraph($ParaIndex, $OldText, $NewText) {
    $r = $d.Paragraphs.Item($ParaIndex).Range
    $found = $r.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, $NewText, 2)
    if (-not $found) {
        Write-Host "WARNING: text not found in paragraph" $ParaIndex "->" $OldText
    }
}

# ---------------------------------------------------------------------------
# Paragraph 6  (body under "Objetivos"): old Objetivos text -> Programa
# resumido text ("1-Processos de Conformação Mecânica" block).
# ---------------------------------------------------------------------------
$old6 = "Apresentar a análise química como ferramenta para o estudo da composição e das propriedades de materiais.^lDesenvolver a competência para formular e compreender problemas relacionados à análise química e buscar de forma autônoma procedimentos adequados para a sua solução. Desenvolver nos alunos a competência técnica para propor experimentos, obter e interpretar resultados analíticos. ^lIncentivar trabalhos em grupo para a solução de problemas, com apresentação de resultados de forma oral e escrita."
$new6 = "1-Processos de Conformação Mecânica^l1.1. Classificação dos Processos de Conformação Mecânica. 1.2. Metalurgia da Conformação Mecânica. 1.3. Mecânica da Conformação. 1.4. Descrição dos Processos de Conformação Plástica: Forjamento, Trefilação, Extrusão, Laminação. 1.5. Processamento de Chapas Metálicas: Estampagem, Dobramento, Calandragem.^l2. Processos de Usinagem ^l2.1 Principais Processos de Usinagem. 2.2. Movimentos e grandezas na usinagem 2.3 Aspectos técnicos do processo. 2.4 Cálculo de potência de usinagem. 2.4. Usinabilidade dos Materiais"
Replace-InParagraph 6 $old6 $new6

# ---------------------------------------------------------------------------
# Paragraph 8 (body under "Docente(s) Responsável(eis)"): the two
# "55840..." runs are replaced, first by the Objetivos text, second by the
# Programa text.
# ---------------------------------------------------------------------------
$old8a = "5840712 - Ângelo Capri Neto"
$new8a = "Apresentar a análise química como ferramenta para o estudo da composição e das propriedades de materiais.^lDesenvolver a competência para formular e compreender problemas relacionados à análise química e buscar de forma autônoma procedimentos adequados para a sua solução. Desenvolver nos alunos a competência técnica para propor experimentos, obter e interpretar resultados analíticos. ^lIncentivar trabalhos em grupo para a solução de problemas, com apresentação de resultados de forma oral e escrita."
Replace-InParagraph 8 $old8a $new8a

$old8b = "5840521 - Rosa Ana Conte"
$new8b = "Introdução à química analítica; Preparação de amostras sólidas e líquidas; Métodos espectroscópicos de análise: interação radiação/matéria, absorção atômica e molecular.^lEspectroscopia UV/Visível: lei de Beer; instrumentação, calibração do equipamento, aplicações e interpretação dos resultados analíticos. Absorção Atômica: instrumentação, calibração do equipamento, identificação e controle de interferências; aplicações e interpretação de resultados analíticos. Emissão Atômica: instrumentação, calibração do equipamento e controle de interferências; aplicações e interpretação de resultados analíticos. Análise de gases em metais: instrumentação e calibração do equipamento; aplicações e interpretação de resultados analíticos."
Replace-InParagraph 8 $old8b $new8b

# ---------------------------------------------------------------------------
# Paragraph 10 (body under "Programa resumido"): replaced by the Avaliação
# "Método" value.
# ---------------------------------------------------------------------------
$old10 = "1-Processos de Conformação Mecânica^l1.1. Classificação dos Processos de Conformação Mecânica. 1.2. Metalurgia da Conformação Mecânica. 1.3. Mecânica da Conformação. 1.4. Descrição dos Processos de Conformação Plástica: Forjamento, Trefilação, Extrusão, Laminação. 1.5. Processamento de Chapas Metálicas: Estampagem, Dobramento, Calandragem.^l2. Processos de Usinagem ^l2.1 Principais Processos de Usinagem. 2.2. Movimentos e grandezas na usinagem 2.3 Aspectos técnicos do processo. 2.4 Cálculo de potência de usinagem. 2.4. Usinabilidade dos Materiais"
$new10 = "A avaliação será feita por meio de duas provas (P1 e P2). A critério do professor, a avaliação poderá ser complementada por meio de trabalhos e/ou relatórios, valendo até 30% da nota das provas."
Replace-InParagraph 10 $old10 $new10

# ---------------------------------------------------------------------------
# Paragraph 12 (body under "Programa"): replaced by the Avaliação
# "Critério" value.
# ---------------------------------------------------------------------------
$old12 = "Introdução à química analítica; Preparação de amostras sólidas e líquidas; Métodos espectroscópicos de análise: interação radiação/matéria, absorção atômica e molecular.^lEspectroscopia UV/Visível: lei de Beer; instrumentação, calibração do equipamento, aplicações e interpretação dos resultados analíticos. Absorção Atômica: instrumentação, calibração do equipamento, identificação e controle de interferências; aplicações e interpretação de resultados analíticos. Emissão Atômica: instrumentação, calibração do equipamento e controle de interferências; aplicações e interpretação de resultados analíticos. Análise de gases em metais: instrumentação e calibração do equipamento; aplicações e interpretação de resultados analíticos."
$new12 = "A nota final (NF) será calculada pela média aritmética das provas. NF=(P1 +P2)/2."
Replace-InParagraph 12 $old12 $new12

# ---------------------------------------------------------------------------
# Paragraph 14 (body under "Avaliação"): three bold labels ("Método:",
# "Critério:", "Norma de recuperação:") stay in place; only the plain-text
# values that follow each label rotate. Use placeholders to avoid the
# later searches accidentally matching text just written by an earlier
# replacement within this same paragraph.
# ---------------------------------------------------------------------------
$oldMetodoValue = "A avaliação será feita por meio de duas provas (P1 e P2). A critério do professor, a avaliação poderá ser complementada por meio de trabalhos e/ou relatórios, valendo até 30% da nota das provas."
$oldCriterioValue = "A nota final (NF) será calculada pela média aritmética das provas. NF=(P1 +P2)/2."
$oldNormaValue = "Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2."

Replace-InParagraph 14 $oldMetodoValue "@@PLACEHOLDER_METODO@@"
Replace-InParagraph 14 $oldCriterioValue "@@PLACEHOLDER_CRITERIO@@"
Replace-InParagraph 14 $oldNormaValue "@@PLACEHOLDER_NORMA@@"

$newMetodoValue = $oldNormaValue
$newCriterioValue = "1. VOGEL, A. L., et al. Análise Química Quantitativa, 6ª Ed., Rio de Janeiro, Livros Técnicos e Científicos Editora S.A., 2003.^l2. SKOOG, D.A. & Jeary, J.J. Principles of Instrumental Analysis, 6th Ed, Saunders College Publishing, 2007.^l3. MITRA, S. Sample Preparation Techniques in Analytical  Chemistry, Wiley & Sons, Hoboken, New Jersey, 2003.^l4. ANDERSON, R. Sample Pretreatment and  separation, Wiley & Sons, New York, 1997"
$newNormaValue = "5840712 - Ângelo Capri Neto"

Replace-InParagraph 14 "@@PLACEHOLDER_METODO@@" $newMetodoValue
Replace-InParagraph 14 "@@PLACEHOLDER_CRITERIO@@" $newCriterioValue
Replace-InParagraph 14 "@@PLACEHOLDER_NORMA@@" $newNormaValue

# ---------------------------------------------------------------------------
# Paragraph 16 (body under "Bibliografia"): replaced by the second
# "Docente(s)" line.
# ---------------------------------------------------------------------------
$old16 = "1. VOGEL, A. L., et al. Análise Química Quantitativa, 6ª Ed., Rio de Janeiro, Livros Técnicos e Científicos Editora S.A., 2003.^l2. SKOOG, D.A. & Jeary, J.J. Principles of Instrumental Analysis, 6th Ed, Saunders College Publishing, 2007.^l3. MITRA, S. Sample Preparation Techniques in Analytical  Chemistry, Wiley & Sons, Hoboken, New Jersey, 2003.^l4. ANDERSON, R. Sample Pretreatment and  separation, Wiley & Sons, New York, 1997"
$new16 = "5840521 - Rosa Ana Conte"
Replace-InParagraph 16 $old16 $new16

Write-Host "Done."
